$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, copying formatting from E1 (bold/border/centered header style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamps for each row (2-18), taken from the target diff
$timestamps = @(
    "2021-10-05 13:42:18.097846",
    "2021-10-05 13:42:18.097859",
    "2021-10-05 13:42:18.097863",
    "2021-10-05 13:42:18.097866",
    "2021-10-05 13:42:18.097870",
    "2021-10-05 13:42:18.097873",
    "2021-10-05 13:42:18.097876",
    "2021-10-05 13:42:18.097879",
    "2021-10-05 13:42:18.097882",
    "2021-10-05 13:42:18.097885",
    "2021-10-05 13:42:18.097888",
    "2021-10-05 13:42:18.097891",
    "2021-10-05 13:42:18.097894",
    "2021-10-05 13:42:18.097897",
    "2021-10-05 13:42:18.097901",
    "2021-10-05 13:42:18.097904",
    "2021-10-05 13:42:18.097907"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
